$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells (Coin, Link, Price, Volume(1h)) to match refreshed crypto data.
# Column D (Price) cells are forced to text format to preserve values such as
# "69.191.45" or "0.0000331" exactly as strings (matching the original inline-string cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.191.45'
$ws.Range("E2").Value = '  +1.38%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.890.32'
$ws.Range("E3").Value = '  -0.21%  '

$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '527.05'
$ws.Range("E5").Value = '  +8.63%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '144.30'
$ws.Range("E6").Value = '  -0.95%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  -1.40%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.717'
$ws.Range("E9").Value = '  -3.10%  '

$ws.Range("E10").Value = '  -6.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000331'
$ws.Range("E11").Value = '  -6.11%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '41.91'
$ws.Range("E12").Value = '  -2.27%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.517.84'
$ws.Range("E13").Value = '  +0.08%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '10.20'
$ws.Range("E14").Value = '  -2.38%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.888.64'
$ws.Range("E15").Value = '  -0.15%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '13.95'
$ws.Range("E16").Value = '  -1.50%  '

$ws.Range("B17").Value = 'TRON'
$ws.Range("C17").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.134'
$ws.Range("E17").Value = '  -1.40%  '

$ws.Range("B18").Value = 'Polygon'
$ws.Range("C18").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.21'
$ws.Range("E18").Value = '  +6.35%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '19.71'
$ws.Range("E19").Value = '  -1.30%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '69.241.12'
$ws.Range("E20").Value = '  +1.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '422.37'
$ws.Range("E21").Value = '  -1.82%  '

$ws.Range("E22").Value = '  -5.20%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.14'
$ws.Range("E23").Value = '  -4.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '87.56'
$ws.Range("E24").Value = '  -1.89%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.99'
$ws.Range("E25").Value = '  +7.92%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.37'
$ws.Range("E26").Value = '  -7.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.55'
$ws.Range("E27").Value = '  -3.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '36.25'
$ws.Range("E28").Value = '  -2.64%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '690.53'
$ws.Range("E29").Value = '  -4.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '13.16'
$ws.Range("E30").Value = '  -1.71%  '

$ws.Range("E31").Value = '  -3.22%  '

$ws.Range("E32").Value = '  -2.49%  '

$ws.Range("E33").Value = '  +10.91%  '

$ws.Range("E34").Value = '  +7.70%  '

$ws.Range("B35").Value = 'NEARProtocol'
$ws.Range("C35").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.90'
$ws.Range("E35").Value = '  -2.74%  '

$ws.Range("B36").Value = 'PEPE'
$ws.Range("C36").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0₃0852'
$ws.Range("E36").Value = '  -3.72%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.82'
$ws.Range("E37").Value = '  -2.13%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.999'
$ws.Range("E38").Value = '  -0.13%  '

$ws.Range("E39").Value = '  +0.51%  '

$ws.Range("B40").Value = 'FirstDigitalUSD'
$ws.Range("C40").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.00%  '

$ws.Range("B41").Value = 'ThetaToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.30'
$ws.Range("E41").Value = '  +7.23%  '

$ws.Range("E42").Value = '  +7.35%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0481'
$ws.Range("E43").Value = '  -2.46%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.76'
$ws.Range("E44").Value = '  -7.03%  '

$ws.Range("E45").Value = '  +1.10%  '

$ws.Range("E46").Value = '  -1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.99'
$ws.Range("E47").Value = '  +6.74%  '

$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.753.06'
$ws.Range("E48").Value = '  +14.82%  '

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '144.89'
$ws.Range("E49").Value = '  +0.74%  '

$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.000268'
$ws.Range("E50").Value = '  +9.13%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.50'
$ws.Range("E51").Value = '  +6.19%  '
